# Zeitplan.xlsx update:
#  - Row 4 ("Woche" 4 / 23.09.2019 - 29.09.2019") now has a "Bemerkung" entry:
#    interview answers have arrived, back on schedule.
#  - Selection moved to the merged "1. Meilenstein" banner row (B6:F6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("F5").Value = "abgeschlossen, Interviewantworten sind angekommen ich bin wieder in Zeitplanung."

$ws.Range("B6:F6").Select() | Out-Null
